# Generate Report for Handback
# Adds a new handback entry (23c49170-0c72-41c0-8ff4-e94a9d827f6c.md) as row 4
# to the "Overview", "zh-cn" and "de-de" tables/worksheets.

$wb = $excel.ActiveWorkbook

$newFileId   = "23c49170-0c72-41c0-8ff4-e94a9d827f6c"
$newFileName = "$newFileId.md"
$newFilePath = "e2e\$newFileName"
$commitHash  = "77beadd421e43a06d591ae9e0631118e91e8dc25"
$srcCommit   = "c9298f28a36f74a09a559650e4f0575ddd8b128c"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $newFileName
$wsOverview.Range("B4").Value = $newFilePath
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e/$newFileName",
    "",
    "",
    $newFilePath
) | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-09-07 06:39:49"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $newFileName
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e/$newFileName",
    "",
    "",
    $newFileName
) | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = "$newFileId.$commitHash.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-07 06:39:36"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $newFileName
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$commitHash/e2e/$newFileName",
    "",
    "",
    $newFileName
) | Out-Null
$wsZh.Range("J4").Value = "$newFileId.$commitHash.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-09-07 06:40:38"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $newFileName
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e/$newFileName",
    "",
    "",
    $newFileName
) | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = "$newFileId.$commitHash.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-07 06:39:49"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $newFileName
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$commitHash/e2e/$newFileName",
    "",
    "",
    $newFileName
) | Out-Null
$wsDe.Range("J4").Value = "$newFileId.$commitHash.de-de.xlf"
$wsDe.Range("K4").Value = "2016-09-07 06:40:56"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = ""

Write-Output "Overview range: $($loOverview.Range.Address())"
Write-Output "zh-cn range: $($loZh.Range.Address())"
Write-Output "de-de range: $($loDe.Range.Address())"
